# Add the two new daily rows (2025-12-19 / serial 46010) for the two
# stations, following the same layout as every prior day in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the previous day's two rows down onto the new
# rows first, so the new cells reuse the existing style/number-format
# indices (date style, "0.00" style, integer style) instead of Excel
# creating brand-new (duplicate) number formats for them.
$ws.Range("A36:F37").Copy()
$ws.Range("A38:F39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 38: 四方坪站 (Sifangping station)
$ws.Range("A38").Value = 46010
$ws.Range("B38").Value = "四方坪站"
$ws.Range("C38").Value = 8259.1
$ws.Range("D38").Value = 6992.14
$ws.Range("E38").Value = 2699.54
$ws.Range("F38").Value = 378

# Row 39: 高岭站 (Gaoling station)
$ws.Range("A39").Value = 46010
$ws.Range("B39").Value = "高岭站"
$ws.Range("C39").Value = 4913.54
$ws.Range("D39").Value = 4016.35
$ws.Range("E39").Value = 1301.6600000000001
$ws.Range("F39").Value = 173

# Scroll the view down a few rows and update the selection to match the
# state the workbook was saved in.
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I37").Select()
